$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before current column B ("Firewall") to hold the
#    new "Region" field. Everything from B onward shifts right by one.
$ws.Columns("B:B").Insert() | Out-Null

# 2) Insert a new row before current row 4 (SIN, the column insert above
#    does not change row numbers) to hold the new SFO site. This pushes
#    SIN/ORD/DUB down by one row.
$ws.Rows("4:4").Insert() | Out-Null

# 3) Append two more rows at the bottom for LHR and AMS.
$ws.Rows("8:8").Insert() | Out-Null
$ws.Rows("9:9").Insert() | Out-Null

# ---- Header row ----
$ws.Range("B1").Value2 = "Region"

# ---- Row 2: GDL ----
$ws.Range("B2").Value2 = "LATAM"

# ---- Row 3: LAX ----
$ws.Range("B3").Value2 = "NA"

# ---- Row 4: SFO (new) ----
$ws.Range("A4").Value2 = "SFO"
$ws.Range("B4").Value2 = "NA"
$ws.Range("C4").Formula = '= CONCATENATE(LOWER(A4),"-fwl01.uwaco.net")'
$ws.Range("D4").Formula = '= CONCATENATE(LOWER(A4),"-iot-loragw01.uwaco.net")'
$ws.Range("E4").Value2 = "10.2.211.50"
$ws.Range("F4").Value2 = "10.2.0.53"
$ws.Range("G4").Value2 = "10.2.0.23"

# ---- Row 5: SIN ----
$ws.Range("B5").Value2 = "APAC"

# ---- Row 6: ORD ----
$ws.Range("B6").Value2 = "NA"

# ---- Row 7: DUB ----
$ws.Range("B7").Value2 = "EMEA"

# ---- Row 8: LHR (new) ----
$ws.Range("A8").Value2 = "LHR"
$ws.Range("B8").Value2 = "EMEA"
$ws.Range("C8").Formula = '= CONCATENATE(LOWER(A8),"-fwl01.uwaco.net")'
$ws.Range("D8").Formula = '= CONCATENATE(LOWER(A8),"-iot-loragw01.uwaco.net")'
$ws.Range("E8").Value2 = "10.7.136.50"
$ws.Range("F8").Value2 = "10.7.0.53"
$ws.Range("G8").Value2 = "10.7.0.23"

# ---- Row 9: AMS (new) ----
$ws.Range("A9").Value2 = "AMS"
$ws.Range("B9").Value2 = "EMEA"
$ws.Range("C9").Formula = '= CONCATENATE(LOWER(A9),"-fwl01.uwaco.net")'
$ws.Range("D9").Formula = '= CONCATENATE(LOWER(A9),"-iot-loragw01.uwaco.net")'
$ws.Range("E9").Value2 = "10.7.225.50"
$ws.Range("F9").Value2 = "10.9.0.53"
$ws.Range("G9").Value2 = "10.9.0.23"

# ---- Sheet view: zoom 150%, new selection ----
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("G13").Select() | Out-Null

Write-Output "done"
